$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 282.72726
$ws.Range("I2").Value = 289.25
$ws.Range("J2").Value = 265.33334
$ws.Range("K2").Value = 289.25
$ws.Range("L2").Value = 265.33334
$ws.Range("M2").Value = -176.25
$ws.Range("N2").Value = -491.33334

$ws.Range("H53").Value = 960.5909
$ws.Range("J53").Value = 937.7
$ws.Range("L53").Value = 937.7
$ws.Range("N53").Value = -2211.7

$ws.Range("H55").Value = 83333610
$ws.Range("J55").Value = 194.5
$ws.Range("L55").Value = 194.5
$ws.Range("N55").Value = -622.5

$ws.Range("H103").Value = 371.75
$ws.Range("I103").Value = 368
$ws.Range("K103").Value = 1104
$ws.Range("M103").Value = -518

$ws.Range("H116").Value = 7500.095
$ws.Range("J116").Value = 6485.7144
$ws.Range("L116").Value = 6485.7144
$ws.Range("N116").Value = -13369.7144

$ws.Range("H138").Value = 3330.55
$ws.Range("J138").Value = 3558.862
$ws.Range("L138").Value = 10676.586
$ws.Range("N138").Value = -20956.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2154.2856
$ws.Range("I2").Value = 1896
$ws.Range("J2").Value = 2498.6667
$ws.Range("K2").Value = 1896
$ws.Range("L2").Value = 2498.6667
$ws.Range("M2").Value = -1783
$ws.Range("N2").Value = -2724.6667

$ws.Range("H4").Value = 410
$ws.Range("I4").Value = 370
$ws.Range("K4").Value = 370
$ws.Range("M4").Value = -254

$ws.Range("H32").Value = 15725.537
$ws.Range("I32").Value = 11547.652
$ws.Range("K32").Value = 11547.652
$ws.Range("M32").Value = -11260.652

$ws.Range("H97").Value = 2966.3235
$ws.Range("I97").Value = 1408
$ws.Range("K97").Value = 1408
$ws.Range("M97").Value = -912

$ws.Range("H116").Value = 2154.2856
$ws.Range("I116").Value = 1896
$ws.Range("J116").Value = 2498.6667
$ws.Range("K116").Value = 1896
$ws.Range("L116").Value = 2498.6667
$ws.Range("M116").Value = 398
$ws.Range("N116").Value = -7086.6667

$ws.Range("H132").Value = 4367.5835
$ws.Range("I132").Value = 3389.6667
$ws.Range("J132").Value = 6323.4165
$ws.Range("K132").Value = 10169.0001
$ws.Range("L132").Value = 18970.2495
$ws.Range("M132").Value = -7639.000100000001
$ws.Range("N132").Value = -24030.2495

$ws.Range("H138").Value = 103333
$ws.Range("J138").Value = 103333
$ws.Range("L138").Value = 103333
$ws.Range("N138").Value = -113613

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2154.2856
$ws.Range("I3").Value = 1896
$ws.Range("J3").Value = 2498.6667
$ws.Range("K3").Value = 1896
$ws.Range("L3").Value = 2498.6667
$ws.Range("M3").Value = -1782
$ws.Range("N3").Value = -2726.6667

$ws.Range("H64").Value = 1578.5
$ws.Range("J64").Value = 1872.75
$ws.Range("L64").Value = 1872.75
$ws.Range("N64").Value = -2322.75

$ws.Range("H67").Value = 1578.5
$ws.Range("J67").Value = 1872.75
$ws.Range("L67").Value = 1872.75
$ws.Range("N67").Value = -3432.75

$ws.Range("H94").Value = 17169.166
$ws.Range("I94").Value = 668.7143
$ws.Range("K94").Value = 668.7143
$ws.Range("M94").Value = -217.7143

$ws.Range("H107").Value = 2131.2666
$ws.Range("I107").Value = 1606.4286
$ws.Range("K107").Value = 1606.4286
$ws.Range("M107").Value = 313.5714

$ws.Range("H134").Value = 5070.466
$ws.Range("I134").Value = 5289
$ws.Range("K134").Value = 15867
$ws.Range("M134").Value = -13332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -213

$ws.Range("H22").Value = 814.0769
$ws.Range("J22").Value = 1012.8
$ws.Range("L22").Value = 1012.8
$ws.Range("N22").Value = -1712.8

$ws.Range("H31").Value = 59360.168
$ws.Range("I31").Value = 3815.3333
$ws.Range("J31").Value = 170449.83
$ws.Range("K31").Value = 3815.3333
$ws.Range("L31").Value = 170449.83
$ws.Range("M31").Value = -3520.3333
$ws.Range("N31").Value = -171039.83

$ws.Range("H34").Value = 59360.168
$ws.Range("I34").Value = 3815.3333
$ws.Range("J34").Value = 170449.83
$ws.Range("K34").Value = 3815.3333
$ws.Range("L34").Value = 170449.83
$ws.Range("M34").Value = -3613.3333
$ws.Range("N34").Value = -170853.83

$ws.Range("H52").Value = 70213.60000000001

$ws.Range("H58").Value = 2605.4666
$ws.Range("I58").Value = 2305.8215
$ws.Range("K58").Value = 2305.8215
$ws.Range("M58").Value = -2102.8215

$ws.Range("H60").Value = 44999.75
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 44999.75
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 44999.75
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -46021.75

$ws.Range("H105").Value = 1429
$ws.Range("I105").Value = 1286.25
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1286.25
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 460.75
$ws.Range("N105").Value = -5494

$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("K113").Value = 500
$ws.Range("M113").Value = 1670

$ws.Range("H132").Value = 3592.375
$ws.Range("I132").Value = 3193.9092
$ws.Range("K132").Value = 9581.7276
$ws.Range("M132").Value = -7051.7276

$ws.Range("H136").Value = 2605.4666
$ws.Range("I136").Value = 2305.8215
$ws.Range("K136").Value = 6917.4645
$ws.Range("M136").Value = -4367.4645

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 438.625
$ws.Range("I86").Value = 418.33334
$ws.Range("J86").Value = 499.5
$ws.Range("K86").Value = 1255.00002
$ws.Range("L86").Value = 1498.5
$ws.Range("M86").Value = -69.00001999999995
$ws.Range("N86").Value = -3870.5

$ws.Range("H89").Value = 438.625
$ws.Range("I89").Value = 418.33334
$ws.Range("J89").Value = 499.5
$ws.Range("K89").Value = 3765.00006
$ws.Range("L89").Value = 4495.5
$ws.Range("M89").Value = 2162.99994
$ws.Range("N89").Value = -16351.5

$ws.Range("H129").Value = 4013.75
$ws.Range("I129").Value = 3237.5
$ws.Range("J129").Value = 4401.875
$ws.Range("K129").Value = 9712.5
$ws.Range("L129").Value = 13205.625
$ws.Range("M129").Value = -4712.5
$ws.Range("N129").Value = -23205.625

$ws.Range("H130").Value = 24981.334
$ws.Range("J130").Value = 24981.334
$ws.Range("L130").Value = 74944.00199999999
$ws.Range("N130").Value = -84984.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3758.6667
$ws.Range("I80").Value = 2999
$ws.Range("J80").Value = 3910.6
$ws.Range("K80").Value = 2999
$ws.Range("L80").Value = 3910.6
$ws.Range("M80").Value = -2001
$ws.Range("N80").Value = -5906.6

$ws.Range("H83").Value = 3758.6667
$ws.Range("I83").Value = 2999
$ws.Range("J83").Value = 3910.6
$ws.Range("K83").Value = 14995
$ws.Range("L83").Value = 19553
$ws.Range("M83").Value = -10003
$ws.Range("N83").Value = -29537

$ws.Range("H113").Value = 6510.3335
$ws.Range("I113").Value = 6142.5713
$ws.Range("K113").Value = 6142.5713
$ws.Range("M113").Value = -3972.5713

$ws.Range("H122").Value = 130401.75
$ws.Range("I122").Value = 252809.5
$ws.Range("K122").Value = 758428.5
$ws.Range("M122").Value = -755978.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 9044.4
$ws.Range("J9").Value = 20999.5
$ws.Range("L9").Value = 20999.5
$ws.Range("N9").Value = -21447.5

$ws.Range("H40").Value = 2494.6086
$ws.Range("I40").Value = 2323.625
$ws.Range("J40").Value = 2885.4285
$ws.Range("K40").Value = 2323.625
$ws.Range("L40").Value = 2885.4285
$ws.Range("M40").Value = -2187.625
$ws.Range("N40").Value = -3157.4285

$ws.Range("H46").Value = 2305.375
$ws.Range("I46").Value = 2057
$ws.Range("K46").Value = 2057
$ws.Range("M46").Value = -1869

$ws.Range("H122").Value = 3851.842
$ws.Range("I122").Value = 3955.8572
$ws.Range("K122").Value = 11867.5716
$ws.Range("M122").Value = -9417.571599999999

$ws.Range("H136").Value = 3917.2856
$ws.Range("J136").Value = 3794
$ws.Range("L136").Value = 11382
$ws.Range("N136").Value = -16482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3057.3408
$ws.Range("I107").Value = 3462.1428
$ws.Range("J107").Value = 2348.9375
$ws.Range("K107").Value = 10386.4284
$ws.Range("L107").Value = 7046.8125
$ws.Range("M107").Value = -8466.428400000001
$ws.Range("N107").Value = -10886.8125

$ws.Range("H136").Value = 8719.362999999999
$ws.Range("I136").Value = 7505.927
$ws.Range("K136").Value = 22517.781
$ws.Range("M136").Value = -19967.781
